$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 43
$ws.Cells.Item(41, 9).Value = 14
$ws.Cells.Item(41, 10).Value = 62.333332
$ws.Cells.Item(41, 11).Value = 14
$ws.Cells.Item(41, 12).Value = 62.333332
$ws.Cells.Item(41, 13).Value = 426
$ws.Cells.Item(41, 14).Value = -942.333332
$ws.Cells.Item(107, 8).Value = 2457.7
$ws.Cells.Item(107, 9).Value = 2014.6666
$ws.Cells.Item(107, 10).Value = 3122.25
$ws.Cells.Item(107, 11).Value = 2014.6666
$ws.Cells.Item(107, 12).Value = 3122.25
$ws.Cells.Item(107, 13).Value = -94.66660000000002
$ws.Cells.Item(107, 14).Value = -6962.25
$ws.Cells.Item(111, 8).Value = 4701.1665
$ws.Cells.Item(111, 9).Value = 4675.684
$ws.Cells.Item(111, 11).Value = 14027.052
$ws.Cells.Item(111, 13).Value = -10960.052
$ws.Cells.Item(115, 8).Value = 688.1667
$ws.Cells.Item(115, 9).Value = 688.1667
$ws.Cells.Item(115, 11).Value = 2064.5001
$ws.Cells.Item(115, 13).Value = -497.5001000000002
$ws.Cells.Item(132, 8).Value = 2878.2593
$ws.Cells.Item(132, 9).Value = 2868.3333
$ws.Cells.Item(132, 11).Value = 8604.999899999999
$ws.Cells.Item(132, 13).Value = -6074.999899999999
$ws.Cells.Item(137, 8).Value = 2245.5
$ws.Cells.Item(137, 9).Value = 2363.5557
$ws.Cells.Item(137, 11).Value = 7090.6671
$ws.Cells.Item(137, 13).Value = -4540.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10531915
$ws.Cells.Item(32, 9).Value = 13700879
$ws.Cells.Item(32, 11).Value = 13700879
$ws.Cells.Item(32, 13).Value = -13700592
$ws.Cells.Item(45, 8).Value = 1772.2222
$ws.Cells.Item(45, 9).Value = 1320
$ws.Cells.Item(45, 11).Value = 1320
$ws.Cells.Item(45, 13).Value = -943
$ws.Cells.Item(61, 8).Value = 41669910
$ws.Cells.Item(61, 9).Value = 55557004
$ws.Cells.Item(61, 10).Value = 8634.833000000001
$ws.Cells.Item(61, 11).Value = 55557004
$ws.Cells.Item(61, 12).Value = 8634.833000000001
$ws.Cells.Item(61, 13).Value = -55556792
$ws.Cells.Item(61, 14).Value = -9058.833000000001
$ws.Cells.Item(63, 8).Value = 4412.147
$ws.Cells.Item(63, 9).Value = 2982.2307
$ws.Cells.Item(63, 10).Value = 9059.375
$ws.Cells.Item(63, 11).Value = 2982.2307
$ws.Cells.Item(63, 12).Value = 9059.375
$ws.Cells.Item(63, 13).Value = -2296.2307
$ws.Cells.Item(63, 14).Value = -10431.375
$ws.Cells.Item(66, 8).Value = 4412.147
$ws.Cells.Item(66, 9).Value = 2982.2307
$ws.Cells.Item(66, 10).Value = 9059.375
$ws.Cells.Item(66, 11).Value = 14911.1535
$ws.Cells.Item(66, 12).Value = 45296.875
$ws.Cells.Item(66, 13).Value = -11479.1535
$ws.Cells.Item(66, 14).Value = -52160.875
$ws.Cells.Item(110, 8).Value = 16198.84
$ws.Cells.Item(110, 9).Value = 17435.727
$ws.Cells.Item(110, 11).Value = 17435.727
$ws.Cells.Item(110, 13).Value = -15390.727
$ws.Cells.Item(122, 8).Value = 3667.375
$ws.Cells.Item(122, 9).Value = 2878.9333
$ws.Cells.Item(122, 10).Value = 4363.0586
$ws.Cells.Item(122, 11).Value = 8636.7999
$ws.Cells.Item(122, 12).Value = 13089.1758
$ws.Cells.Item(122, 13).Value = -6186.7999
$ws.Cells.Item(122, 14).Value = -17989.1758
$ws.Cells.Item(136, 8).Value = 41669910
$ws.Cells.Item(136, 9).Value = 55557004
$ws.Cells.Item(136, 10).Value = 8634.833000000001
$ws.Cells.Item(136, 11).Value = 166671012
$ws.Cells.Item(136, 12).Value = 25904.499
$ws.Cells.Item(136, 13).Value = -166668462
$ws.Cells.Item(136, 14).Value = -31004.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2038.6154
$ws.Cells.Item(94, 9).Value = 1883.6666
$ws.Cells.Item(94, 10).Value = 2249.9092
$ws.Cells.Item(94, 11).Value = 1883.6666
$ws.Cells.Item(94, 12).Value = 2249.9092
$ws.Cells.Item(94, 13).Value = -1432.6666
$ws.Cells.Item(94, 14).Value = -3151.9092
$ws.Cells.Item(105, 8).Value = 34120
$ws.Cells.Item(105, 9).Value = 50505
$ws.Cells.Item(105, 11).Value = 50505
$ws.Cells.Item(105, 13).Value = -48758
$ws.Cells.Item(134, 8).Value = 1805.9608
$ws.Cells.Item(134, 9).Value = 1833.88
$ws.Cells.Item(134, 10).Value = 410
$ws.Cells.Item(134, 11).Value = 5501.64
$ws.Cells.Item(134, 12).Value = 1230
$ws.Cells.Item(134, 13).Value = -2966.64
$ws.Cells.Item(134, 14).Value = -6300

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 6029.8
$ws.Cells.Item(62, 9).Value = 3999.5
$ws.Cells.Item(62, 10).Value = 7383.3335
$ws.Cells.Item(62, 11).Value = 3999.5
$ws.Cells.Item(62, 12).Value = 7383.3335
$ws.Cells.Item(62, 13).Value = -3375.5
$ws.Cells.Item(62, 14).Value = -8631.333500000001
$ws.Cells.Item(65, 8).Value = 6029.8
$ws.Cells.Item(65, 9).Value = 3999.5
$ws.Cells.Item(65, 10).Value = 7383.3335
$ws.Cells.Item(65, 11).Value = 19997.5
$ws.Cells.Item(65, 12).Value = 36916.6675
$ws.Cells.Item(65, 13).Value = -16877.5
$ws.Cells.Item(65, 14).Value = -43156.6675
$ws.Cells.Item(99, 8).Value = 16121
$ws.Cells.Item(99, 9).Value = 19297.166
$ws.Cells.Item(99, 11).Value = 19297.166
$ws.Cells.Item(99, 13).Value = -17799.166
$ws.Cells.Item(107, 8).Value = 1213.5714
$ws.Cells.Item(107, 9).Value = 1448
$ws.Cells.Item(107, 10).Value = 1037.75
$ws.Cells.Item(107, 11).Value = 1448
$ws.Cells.Item(107, 12).Value = 1037.75
$ws.Cells.Item(107, 13).Value = 472
$ws.Cells.Item(107, 14).Value = -4877.75
$ws.Cells.Item(112, 8).Value = 55133.332
$ws.Cells.Item(112, 9).Value = 25000
$ws.Cells.Item(112, 10).Value = 70200
$ws.Cells.Item(112, 11).Value = 25000
$ws.Cells.Item(112, 12).Value = 70200
$ws.Cells.Item(112, 13).Value = -23523
$ws.Cells.Item(112, 14).Value = -73154
$ws.Cells.Item(126, 8).Value = 16121
$ws.Cells.Item(126, 9).Value = 19297.166
$ws.Cells.Item(126, 11).Value = 57891.49800000001
$ws.Cells.Item(126, 13).Value = -55421.49800000001
$ws.Cells.Item(134, 8).Value = 1290.8077
$ws.Cells.Item(134, 9).Value = 1148.375
$ws.Cells.Item(134, 11).Value = 3445.125
$ws.Cells.Item(134, 13).Value = -910.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 2820
$ws.Cells.Item(3, 9).Value = 2820
$ws.Cells.Item(3, 11).Value = 8460
$ws.Cells.Item(3, 13).Value = -8348
$ws.Cells.Item(98, 8).Value = 2169.8
$ws.Cells.Item(98, 10).Value = 2359.3845
$ws.Cells.Item(98, 12).Value = 7078.1535
$ws.Cells.Item(98, 14).Value = -10074.1535
$ws.Cells.Item(107, 8).Value = 1299.8
$ws.Cells.Item(107, 9).Value = 300
$ws.Cells.Item(107, 10).Value = 1549.75
$ws.Cells.Item(107, 11).Value = 900
$ws.Cells.Item(107, 12).Value = 4649.25
$ws.Cells.Item(107, 13).Value = 1020
$ws.Cells.Item(107, 14).Value = -8489.25
$ws.Cells.Item(129, 8).Value = 4576.731
$ws.Cells.Item(129, 9).Value = 4449.875
$ws.Cells.Item(129, 10).Value = 4633.1113
$ws.Cells.Item(129, 11).Value = 13349.625
$ws.Cells.Item(129, 12).Value = 13899.3339
$ws.Cells.Item(129, 13).Value = -8349.625
$ws.Cells.Item(129, 14).Value = -23899.3339
$ws.Cells.Item(132, 8).Value = 1756987.5
$ws.Cells.Item(132, 10).Value = 3924684
$ws.Cells.Item(132, 12).Value = 35322156
$ws.Cells.Item(132, 14).Value = -35327216

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1275.8572
$ws.Cells.Item(97, 9).Value = 415.58334
$ws.Cells.Item(97, 11).Value = 415.58334
$ws.Cells.Item(97, 13).Value = 80.41665999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3754.9524
$ws.Cells.Item(40, 9).Value = 3110.818
$ws.Cells.Item(40, 11).Value = 3110.818
$ws.Cells.Item(40, 13).Value = -2974.818
$ws.Cells.Item(68, 8).Value = 2778.1052
$ws.Cells.Item(68, 9).Value = 2516.7646
$ws.Cells.Item(68, 10).Value = 4999.5
$ws.Cells.Item(68, 11).Value = 2516.7646
$ws.Cells.Item(68, 12).Value = 4999.5
$ws.Cells.Item(68, 13).Value = -1767.7646
$ws.Cells.Item(68, 14).Value = -6497.5
$ws.Cells.Item(71, 8).Value = 2778.1052
$ws.Cells.Item(71, 9).Value = 2516.7646
$ws.Cells.Item(71, 10).Value = 4999.5
$ws.Cells.Item(71, 11).Value = 12583.823
$ws.Cells.Item(71, 12).Value = 24997.5
$ws.Cells.Item(71, 13).Value = -8839.823
$ws.Cells.Item(71, 14).Value = -32485.5
$ws.Cells.Item(100, 8).Value = 4434.3335
$ws.Cells.Item(100, 10).Value = 5151.5
$ws.Cells.Item(100, 12).Value = 5151.5
$ws.Cells.Item(100, 14).Value = -6233.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 7913.643
$ws.Cells.Item(62, 9).Value = 8333.333000000001
$ws.Cells.Item(62, 10).Value = 7799.1816
$ws.Cells.Item(62, 11).Value = 8333.333000000001
$ws.Cells.Item(62, 12).Value = 7799.1816
$ws.Cells.Item(62, 13).Value = -7709.333000000001
$ws.Cells.Item(62, 14).Value = -9047.1816
$ws.Cells.Item(65, 8).Value = 7913.643
$ws.Cells.Item(65, 9).Value = 8333.333000000001
$ws.Cells.Item(65, 10).Value = 7799.1816
$ws.Cells.Item(65, 11).Value = 41666.665
$ws.Cells.Item(65, 12).Value = 38995.908
$ws.Cells.Item(65, 13).Value = -38546.665
$ws.Cells.Item(65, 14).Value = -45235.908
$ws.Cells.Item(96, 8).Value = 6825
$ws.Cells.Item(96, 9).Value = 5567.2
$ws.Cells.Item(96, 10).Value = 9969.5
$ws.Cells.Item(96, 11).Value = 5567.2
$ws.Cells.Item(96, 12).Value = 9969.5
$ws.Cells.Item(96, 13).Value = -4194.2
$ws.Cells.Item(96, 14).Value = -12715.5
$ws.Cells.Item(112, 8).Value = 23849.4
$ws.Cells.Item(112, 10).Value = 23849.4
$ws.Cells.Item(112, 12).Value = 23849.4
$ws.Cells.Item(112, 14).Value = -26803.4
$ws.Cells.Item(132, 8).Value = 4389.531
$ws.Cells.Item(132, 9).Value = 4499.5
$ws.Cells.Item(132, 10).Value = 3900.7778
$ws.Cells.Item(132, 11).Value = 13498.5
$ws.Cells.Item(132, 12).Value = 11702.3334
$ws.Cells.Item(132, 13).Value = -10968.5
$ws.Cells.Item(132, 14).Value = -16762.3334
